# Apply crypto price/volume updates from the automated data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Force the cell to stay a text value (matches original inlineStr cells)
    # rather than letting Excel auto-convert numeric-looking strings
    # (e.g. "1.00" or "7.50") into numbers and lose formatting/trailing zeros.
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextCell "D2" "58.868.94"
$ws.Range("E2").Value = "  +2.50%  "
Set-TextCell "D3" "2.518.56"
$ws.Range("E3").Value = "  +3.79%  "
$ws.Range("E4").Value = "  +0.16%  "
Set-TextCell "D5" "532.78"
$ws.Range("E5").Value = "  +4.53%  "
Set-TextCell "D6" "135.06"
$ws.Range("E6").Value = "  +5.43%  "
Set-TextCell "D7" "1.00"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  +3.48%  "
Set-TextCell "D9" "2.520.25"
$ws.Range("E9").Value = "  +3.41%  "
Set-TextCell "D10" "0.0993"
$ws.Range("E10").Value = "  +5.18%  "
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("E12").Value = "  +2.20%  "
Set-TextCell "D13" "0.334"
$ws.Range("E13").Value = "  +1.50%  "
Set-TextCell "D14" "2.957.82"
$ws.Range("E14").Value = "  +3.43%  "
Set-TextCell "D15" "58.925.75"
$ws.Range("E15").Value = "  +2.74%  "
Set-TextCell "D16" "22.43"
$ws.Range("E16").Value = "  +3.75%  "
$ws.Range("E17").Value = "  +3.87%  "
Set-TextCell "D18" "2.507.24"
$ws.Range("E18").Value = "  +3.02%  "
Set-TextCell "D19" "10.69"
$ws.Range("E19").Value = "  +2.99%  "
Set-TextCell "D20" "4.24"
$ws.Range("E20").Value = "  +3.83%  "
Set-TextCell "D21" "322.57"
$ws.Range("E21").Value = "  +2.75%  "
Set-TextCell "D22" "6.14"
$ws.Range("E22").Value = "  +9.38%  "
Set-TextCell "D23" "1.00"
$ws.Range("E23").Value = "  -0.04%  "
Set-TextCell "D24" "65.84"
$ws.Range("E24").Value = "  +4.24%  "
Set-TextCell "D25" "0.410"
$ws.Range("E25").Value = "  +2.01%  "
Set-TextCell "D26" "0.995"
$ws.Range("E26").Value = "  -0.14%  "
Set-TextCell "D27" "0.161"
$ws.Range("E27").Value = "  +1.43%  "
Set-TextCell "D28" "7.50"
$ws.Range("E28").Value = "  +4.19%  "
Set-TextCell "D29" "0.0₃0765"
$ws.Range("E29").Value = "  +7.17%  "
Set-TextCell "D30" "171.67"
$ws.Range("E30").Value = "  +1.09%  "
Set-TextCell "D31" "1.74"
$ws.Range("E31").Value = "  +5.10%  "
$ws.Range("E32").Value = "  +4.91%  "
Set-TextCell "D33" "6.30"
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("E34").Value = "  +0.02%  "
Set-TextCell "D35" "0.997"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  +3.50%  "
$ws.Range("E37").Value = "  -0.96%  "
Set-TextCell "D38" "3.98"
$ws.Range("E38").Value = "  +2.26%  "
Set-TextCell "D39" "1.52"
$ws.Range("E39").Value = "  +5.45%  "
$ws.Range("E40").Value = "  +1.40%  "
Set-TextCell "D41" "0.788"
$ws.Range("E41").Value = "  +3.01%  "
Set-TextCell "D42" "5.20"
$ws.Range("E42").Value = "  +7.30%  "
Set-TextCell "D43" "279.49"
$ws.Range("E43").Value = "  +3.22%  "
Set-TextCell "D44" "3.49"
$ws.Range("E44").Value = "  +3.96%  "
Set-TextCell "D45" "132.08"
$ws.Range("E45").Value = "  +10.29%  "
$ws.Range("E46").Value = "  +3.13%  "
$ws.Range("E47").Value = "  +3.07%  "
$ws.Range("E48").Value = "  +6.04%  "
Set-TextCell "D49" "0.0218"
$ws.Range("E49").Value = "  +4.72%  "
Set-TextCell "D50" "17.13"
$ws.Range("E50").Value = "  +4.22%  "
Set-TextCell "D51" "1.759.20"
$ws.Range("E51").Value = "  +3.67%  "
